$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    ,@(-18.77383784294282, 2.232647813278923, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, 2.49815683730869, -18.77383784294282, -18.77383784294282)
    ,@(-18.77383784294282, 2.157428450105745, 2.886148035767876, -18.77383784294282, 2.591726898477692, -18.77383784294282, 1.880320030434071, -18.77383784294282, 2.446833922759293, -18.77383784294282)
    ,@(-18.77383784294282, 1.263758709243291, -18.77383784294282, -18.77383784294282, -18.77383784294282, 2.110884941596344, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282)
    ,@(2.888384158256683, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282)
    ,@(3.654842564353805, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, 1.54328435293061, -18.77383784294282, 2.22477806489103)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, 2.547066798616085, -18.77383784294282, -18.77383784294282, -18.77383784294282, 1.361872695500483)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, 4.321925037078887, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, 2.270765058669093, 1.610685604070312)
    ,@(-18.77383784294282, -18.77383784294282, 1.717409502706474, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, 2.059086929402126)
    ,@(-18.77383784294282, -18.77383784294282, -0.9688434702191976, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, 2.322659207865831, -18.77383784294282)
    ,@(-18.77383784294282, 0.5987913217339964, -0.1832179012188619, -18.77383784294282, -18.77383784294282, -18.77383784294282, 0.4886355610784074, 0.8835870321797273, 1.140256723691478, -18.77383784294282)
    ,@(-18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, -18.77383784294282, 0.3746412218445586, 0.8171091487461212, 1.330705885121106, -18.77383784294282)
    ,@(-18.77383784294282, -18.77383784294282, 1.555315761870141, -18.77383784294282, -18.77383784294282, -18.77383784294282, 1.847746359354018, 2.06734333017109, -18.77383784294282, -18.77383784294282)
    ,@(-18.77383784294282, 1.65465421470698, 2.318445077107227, -18.77383784294282, 3.80444201034333, -18.77383784294282, 2.177433352387251, 1.863063748516958, -18.77383784294282, 2.468373897467063)
    ,@(-18.77383784294282, 1.916358825130363, -18.77383784294282, -18.77383784294282, -18.77383784294282, 3.298064728032957, 2.458358793416735, -18.77383784294282, -18.77383784294282, -18.77383784294282)
)

for ($r = 0; $r -lt $values.Count; $r++) {
    for ($c = 0; $c -lt $values[$r].Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $values[$r][$c]
    }
}
